$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 24 - this shifts the existing data rows 24:95
# down to 25:96 (and updates the used-range dimension accordingly).
$ws.Rows("24:24").Insert()

# Populate the newly inserted row 24 with the new record.
$ws.Range("A24").Value = 11
$ws.Range("B24").Value = "Vega Monumental Concepción"
$ws.Range("C24").Value = "Bíobío"
$ws.Range("D24").Value = "2022-04-08"
$ws.Range("E24").Value = 8
$ws.Range("F24").Value = 100112021
$ws.Range("G24").Value = "Ají"
$ws.Range("H24").Value = "Chilena(o)"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 25
$ws.Range("K24").Value = 20000
$ws.Range("L24").Value = 21000
$ws.Range("M24").Value = 20600
$ws.Range("N24").Value = '$/saco 25 kilos'
$ws.Range("O24").Value = "Región Metropolitana"
$ws.Range("P24").Value = 824
$ws.Range("Q24").Value = 25
$ws.Range("R24").Value = "Hortaliza"
